# PEAJES.xlsx update — "se actualizo los campos"
# Replaces the toll-invoice rows (dates, invoice numbers, amounts) with the
# new data set, drops the obsolete RUC column, marks the invoice-number
# column as text, renames the sheet and trims the used range to 12 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the sheet ---------------------------------------------------
$ws.Name = "set2023"

# --- new data (row -> A date, B invoice-series, C invoice-number,
#               D rate, E document-id, G company) -----------------------
$data = @(
    @{ Row=1;  A="28.08.2023"; B="F260"; C=367658;  D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=2;  A="19.09.2023"; B="F751"; C=2024734; D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=3;  A="29.08.2023"; B="F158"; C=344693;  D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=4;  A="05.09.2023"; B="F261"; C=141683;  D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=5;  A="05.09.2023"; B="F651"; C=2232357; D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=6;  A="06.09.2023"; B="F255"; C=3399744; D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=7;  A="14.09.2023"; B="F159"; C=2081760; D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=8;  A="18.09.2023"; B="F159"; C=2087717; D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=9;  A="18.09.2023"; B="F152"; C=3095548; D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=10; A="19.09.2023"; B="F257"; C=1043782; D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=11; A="19.09.2023"; B="F160"; C=1172087; D=6.6;    E=20523621212; G="LIMA EXPRESA S.A.C." },
    @{ Row=12; A="19.09.2023"; B="FA17"; C=155632;  D=226.56; E=20608300393; G="COMPAÑIA FOOD RETAIL S.A.C." }
)

# Column A holds day.month.year text like "05.09.2023" which Excel's smart
# entry would otherwise silently re-interpret as a real date for any value
# whose day number is <= 12. Force text mode first, then strip the style
# back off (the authored file keeps these as plain shared-string cells with
# no explicit style), so only the literal text is stored.
$ws.Range("A1:A12").NumberFormat = "@"

foreach ($r in $data) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A   # A - date (text)
    $ws.Cells.Item($row, 2).Value = $r.B   # B - invoice series
    $ws.Cells.Item($row, 3).Value = $r.C   # C - invoice number (numeric, text-formatted below)
    $ws.Cells.Item($row, 4).Value = $r.D   # D - rate
    $ws.Cells.Item($row, 5).Value = $r.E   # E - document id
    $ws.Cells.Item($row, 7).Value = $r.G   # G - company name
}

$ws.Range("A1:A12").Style = "Normal"

# Invoice numbers in column C stay numeric but are displayed/stored with a
# text number format (numFmtId 49 / "@"), matching the authored style.
$ws.Range("C1:C12").NumberFormat = "@"

# Column F (old invoice text) is gone, and the obsolete RUC values that used
# to sit in column G are gone too (G now only carries the company name).
$ws.Range("F1:F15").ClearContents() | Out-Null

# Only 12 data rows remain.
$ws.Range("A13:G15").ClearContents() | Out-Null

# --- selection matches the authored file --------------------------------
$ws.Range("C1:C12").Select() | Out-Null
